# iResearch (Pie Chart) block: fold the standalone chart picture + the
# second ("Source: ...") table into the first ("Chart Title") table so
# the whole placeholder lives in a single 3-row table, exactly like the
# companion bar-chart block. Also shrinks the source row down to a
# trHeight of 189 twips (~0.01in after rounding) and refreshes the
# table/cell borders+shading to the current theme-based style.

$d = $word.ActiveDocument

# Sanity-check the shape we expect to be editing before we blow it away.
if ($d.Tables.Count -ne 2) {
    throw "expected 2 tables (title table + source table), found $($d.Tables.Count)"
}

# NB: kept on one line (no pretty-printed indentation) on purpose --
# InsertXML preserves whitespace verbatim inside foreign content like
# <w:drawing>/<wp:inline>, so indenting this heredoc would leak stray
# whitespace text nodes into the drawing markup.
$newBodyXml = '<w:tbl><w:tblPr><w:tblW w:w="10350" w:type="dxa"/><w:tblInd w:w="108" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/><w:tblLayout w:type="fixed"/><w:tblLook w:val="01E0"/></w:tblPr><w:tblGrid><w:gridCol w:w="10350"/></w:tblGrid><w:tr><w:trPr><w:trHeight w:val="144"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="10350" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="396295" w:themeColor="accent6" w:themeShade="80"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="396295" w:themeColor="accent6" w:themeShade="80"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="iRPlaceholderTitle"/></w:pPr><w:r><w:t>Chart Title</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="3183"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="10350" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="396295" w:themeColor="accent6" w:themeShade="80"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="396295" w:themeColor="accent6" w:themeShade="80"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0"><wp:extent cx="3083139" cy="1971304"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="1" name="Chart 1"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/chart"><c:chart xmlns:c="http://schemas.openxmlformats.org/drawingml/2006/chart" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" r:id="rId8"/></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="189"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="10350" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="396295" w:themeColor="accent6" w:themeShade="80"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="iRPlaceholderNote"/></w:pPr><w:r><w:t>Source: HC</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p/>'

# $d.Content spans the whole body story except the final sectPr, i.e.
# exactly the two tables + the drawing paragraph + the blank spacer
# paragraphs that sit between/after them. InsertXML replaces the
# targeted range's contents wholesale, which is exactly what we want:
# swap all of that for the single merged table (plus the trailing blank
# paragraph Word always keeps before the section break).
$d.Content.InsertXML($newBodyXml)

Write-Output "Tables after edit: $($d.Tables.Count)"
Write-Output "Rows in table 1: $($d.Tables.Item(1).Rows.Count)"
